# Update want-to-go counts (column F) and minimum ticket price (column G)
# for several events on the "展览" (Exhibitions) and "全部类型" (All Types) sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 1903

$ws1.Range("F5").Value = 181
$ws1.Range("G5").Value = 50

$ws1.Range("F6").Value = 2743
$ws1.Range("F7").Value = 191
$ws1.Range("F8").Value = 99
$ws1.Range("F9").Value = 181
$ws1.Range("F10").Value = 1577
$ws1.Range("F11").Value = 557

$ws1.Range("F17").Value = 10
$ws1.Range("F18").Value = 218

$ws1.Range("F21").Value = 15
$ws1.Range("F22").Value = 11
$ws1.Range("F23").Value = 222
$ws1.Range("F24").Value = 69
$ws1.Range("F25").Value = 1752

$ws1.Range("F27").Value = 420
$ws1.Range("F28").Value = 83
$ws1.Range("F29").Value = 571

$ws1.Range("F31").Value = 311
$ws1.Range("F32").Value = 451

# ---- Sheet: 全部类型 (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 1903

$ws4.Range("F6").Value = 181
$ws4.Range("G6").Value = 50

$ws4.Range("F7").Value = 2743
$ws4.Range("F8").Value = 191
$ws4.Range("F9").Value = 99
$ws4.Range("F10").Value = 181
$ws4.Range("F11").Value = 1577
$ws4.Range("F12").Value = 557

$ws4.Range("F18").Value = 10
$ws4.Range("F19").Value = 218

$ws4.Range("F22").Value = 15
$ws4.Range("F23").Value = 11
$ws4.Range("F24").Value = 222
$ws4.Range("F25").Value = 69
$ws4.Range("F26").Value = 1752

$ws4.Range("F28").Value = 420
$ws4.Range("F29").Value = 83
$ws4.Range("F30").Value = 571

$ws4.Range("F32").Value = 311
$ws4.Range("F33").Value = 451
